# "add area to Q files stn5"
# Adds cross-sectional Area (column G/H) alongside the existing discharge
# (Q) computation, and mirrors the two totals (Atotal, Qtotal) side by
# side in columns J/K for a quick-glance summary.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 headers -----------------------------------------------------
# D1 keeps its text ("segment") - only its shared-string slot shifts
# because the old "-" placeholder string is no longer used once B2/C2
# become real numbers (see below).
$ws.Range("D1").Value = "segment"
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# --- B2 / C2: first station's velocity/depth are real 0s, not "-" -----
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0

# --- Column G: cross-sectional Area per segment ------------------------
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# --- Column H: Atotal (sum of segment areas) ----------------------------
$ws.Range("H2").Formula = "=SUM(G2:G11)"

# --- Columns J/K: side-by-side summary of the two totals ---------------
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# --- Selection left on the new summary cells, matching the saved file --
$ws.Range("J2:K2").Select() | Out-Null

$wb.Application.Calculate() | Out-Null
